# Auto-generated script to update cryptos.xlsx price/volume figures
# (recreates the Sat Apr  8 07:27:42 UTC 2023 GitHub Actions cryptos-list refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.256.24"
$ws.Range("E2").Value = "  +0.81%  "
$ws.Range("D3").Value = "1.884.44"
$ws.Range("E3").Value = "  +0.79%  "
$ws.Range("D4").Value = "'1.005"
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "'313.95"
$ws.Range("E5").Value = "  +0.44%  "
$ws.Range("D6").Value = "'1.005"
$ws.Range("E6").Value = "  +0.19%  "
$ws.Range("D7").Value = "'0.5141"
$ws.Range("E7").Value = "  +0.78%  "
$ws.Range("D8").Value = "'0.3908"
$ws.Range("E8").Value = "  +2.81%  "
$ws.Range("D9").Value = "'0.08366"
$ws.Range("E9").Value = "  +0.80%  "
$ws.Range("E10").Value = "  +0.92%  "
$ws.Range("D11").Value = "'41.65"
$ws.Range("E11").Value = "  +0.46%  "
$ws.Range("D12").Value = "'6.242"
$ws.Range("E12").Value = "  +0.50%  "
$ws.Range("D13").Value = "'20.77"
$ws.Range("E13").Value = "  +1.51%  "
$ws.Range("D14").Value = "1.880.16"
$ws.Range("E14").Value = "  +0.57%  "
$ws.Range("D15").Value = "'7.297"
$ws.Range("E15").Value = "  +1.67%  "
$ws.Range("D16").Value = "'1.005"
$ws.Range("E16").Value = "  +0.11%  "
$ws.Range("E17").Value = "  +1.24%  "
$ws.Range("D18").Value = "'91.46"
$ws.Range("E18").Value = "  +0.86%  "
$ws.Range("D19").Value = "'0.06666"
$ws.Range("E19").Value = "  +0.58%  "
$ws.Range("E20").Value = "  +0.07%  "
$ws.Range("D21").Value = "'1.005"
$ws.Range("E21").Value = "  +0.29%  "
$ws.Range("E22").Value = "  +1.23%  "
$ws.Range("D23").Value = "28.292.03"
$ws.Range("E23").Value = "  +0.82%  "
$ws.Range("D24").Value = "'11.18"
$ws.Range("E24").Value = "  +0.69%  "
$ws.Range("D25").Value = "'2.269"
$ws.Range("E25").Value = "  +0.47%  "
$ws.Range("D26").Value = "2.093.17"
$ws.Range("E26").Value = "  +0.41%  "
$ws.Range("E27").Value = "  -1.90%  "
$ws.Range("D28").Value = "'158.90"
$ws.Range("E28").Value = "  +1.15%  "
$ws.Range("D29").Value = "'20.68"
$ws.Range("E29").Value = "  +0.96%  "
$ws.Range("D30").Value = "'125.54"
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("D31").Value = "'0.1066"
$ws.Range("E31").Value = "  +0.92%  "
$ws.Range("D32").Value = "'1.046"
$ws.Range("E32").Value = "  +0.57%  "
$ws.Range("D33").Value = "'5.893"
$ws.Range("E33").Value = "  +5.44%  "
$ws.Range("D34").Value = "'3.599"
$ws.Range("E34").Value = "  -0.22%  "
$ws.Range("D35").Value = "'9.803"
$ws.Range("E35").Value = "  +1.38%  "
$ws.Range("D36").Value = "'0.02458"
$ws.Range("E36").Value = "  +1.16%  "
$ws.Range("D37").Value = "'0.06571"
$ws.Range("E37").Value = "  +0.55%  "
$ws.Range("D38").Value = "'0.2195"
$ws.Range("E38").Value = "  +1.63%  "
$ws.Range("D39").Value = "'1.213"
$ws.Range("E39").Value = "  +0.62%  "
$ws.Range("D40").Value = "'0.6548"
$ws.Range("E40").Value = "  +2.01%  "
$ws.Range("D41").Value = "'5.040"
$ws.Range("E41").Value = "  +3.50%  "
$ws.Range("E42").Value = "  -0.47%  "
$ws.Range("D43").Value = "'11.31"
$ws.Range("E43").Value = "  +0.48%  "
$ws.Range("D44").Value = "'0.6140"
$ws.Range("E44").Value = "  +0.65%  "
$ws.Range("D45").Value = "'13.15"
$ws.Range("E45").Value = "  +1.51%  "
$ws.Range("D46").Value = "'1.294"
$ws.Range("E46").Value = "  +0.74%  "
$ws.Range("D47").Value = "'3.681"
$ws.Range("E47").Value = "  +0.53%  "
$ws.Range("D48").Value = "'2.019"
$ws.Range("E48").Value = "  +1.17%  "
$ws.Range("D49").Value = "'1.237"
$ws.Range("D50").Value = "'121.72"
$ws.Range("E50").Value = "  +0.45%  "
$ws.Range("D51").Value = "'78.91"
$ws.Range("E51").Value = "  -1.22%  "
